$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.275139689445496
$ws.Range("B1").Value = 1.483788728713989
$ws.Range("C1").Value = 1.502338528633118
$ws.Range("D1").Value = 1.565640807151794
$ws.Range("E1").Value = 1.257529377937317
